# Applies the edits described in the diff:
#  1. "...具体实现）等。" -> "...具体实现），划分训练集，验证集和测试集等。"
#  2. "网络训练要考虑的内容..." -> "网络训练和测试要考虑的内容..."
#  3. "每一个部分的训练选择需要" -> "每一个部分的训练选择和测试结果需要"
#  4. "模型文档" (bold, red) -> "模型训练与测试文档（Model_TrainTest）" (bold, red)
#  5. "问题解决文档" (bold, red) -> "问题解决文档（Problem_Handler）" (bold, red)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: insert "，划分训练集，验证集和测试集" right before the "等" that
# follows "...（解决数据不均衡具体实现）"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("具体实现）等。", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins = $rng.Duplicate
$ins.MoveStart(1, 5)    # move past "具体实现）" (5 characters) to just before "等"
$ins.Collapse(1)
$ins.InsertBefore("，划分训练集，验证集和测试集")

# ---------------------------------------------------------------------
# Edit 2: "网络训练要考虑的内容..." -> "网络训练和测试要考虑的内容..."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("网络训练要考虑的内容包括但不限于", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins = $rng.Duplicate
$ins.MoveStart(1, 4)    # move past "网络训练" (4 characters)
$ins.Collapse(1)
$ins.InsertBefore("和测试")

# ---------------------------------------------------------------------
# Edit 3: "每一个部分的训练选择需要" -> "每一个部分的训练选择和测试结果需要"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("每一个部分的训练选择需要", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins = $rng.Duplicate
$ins.MoveStart(1, 10)   # move past "每一个部分的训练选择" (10 characters)
$ins.Collapse(1)
$ins.InsertBefore("和测试结果")

# ---------------------------------------------------------------------
# Edit 4: "模型文档" (bold, red) -> "模型训练与测试文档（Model_TrainTest）" (bold, red)
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("模型文档", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins = $rng.Duplicate
$ins.MoveStart(1, 2)    # move past "模型" (2 characters), leave "文档" untouched after
$ins.Collapse(1)
$ins.InsertBefore("训练与测试")

$ins2 = $d.Content
$ins2.Find.Execute("模型训练与测试文档", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins2.Collapse(0)
$ins2.InsertBefore("（Model_TrainTest）")

# ---------------------------------------------------------------------
# Edit 5: "问题解决文档" (bold, red) -> append "（Problem_Handler）" (bold, red)
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("问题解决文档", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins = $rng.Duplicate
$ins.Collapse(0)
$ins.InsertBefore("（Problem_Handler）")

Write-Output "done"
